$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "recherche information moteur/ achat pièce moteur/batterie"
$ws.Range("D3").Value = "recherche informations transmission, cerveau moteur, achat pièces complemetaires"

$ws.Range("C4").Value = "développement pièces transmission"
$ws.Range("D4").Value = "développement pièces direction"

$ws.Range("C5").Value = "imprimer pièces faites au td d'avant "
$ws.Range("D5").Value = "imprimer pièces faites au td d'avant "

$ws.Range("C7").Value = "gestion direction avec télécomande (code)"
$ws.Range("D7").Value = "gestion vitesse avec télécommande"

$ws.Range("C8").Value = "gestion direction avec télécomande (code)"
$ws.Range("D8").Value = "gestion vitesse avec télécommande"

$ws.Range("C9").Value = "création télécommande"
$ws.Range("D9").Value = "assemblage voiture"

$ws.Range("D13").Select()
